$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "2020" year column (L) should actually be "2021", and a brand new
# "2022" year column (M) is appended after it, carrying the same
# formatting as the existing year columns.

# Fix the mislabeled year in L3 (was 2020, should be 2021).
$ws.Range("L3").Value = 2021

# Add the new 2022 column: clone L3/L4 formatting into M3/M4, then set
# the new values (header year + repeated metric value).
$ws.Range("L3").Copy($ws.Range("M3"))
$ws.Range("M3").Value = 2022

$ws.Range("L4").Copy($ws.Range("M4"))
$ws.Range("M4").Value = 6.18

# Move the active selection to M9, matching the saved view state.
$ws.Range("M9").Select() | Out-Null
